$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that lived on D2 (richfey@gmail.com)
$ws.Range("D2").Hyperlinks.Delete()

# Clear all existing header/data content (B1:L2) - this keeps D2's style
# (the Hyperlink cell style) while dropping its value, matching the
# "emptied but still styled" target cell.
$ws.Range("A1:L2").ClearContents()

# Turn sheet off right-to-left display
$excel.ActiveWindow.DisplayRightToLeft = $false

# Re-write the header row, shifted one column to the left (A1:K1).
# Column order: FULL NAME, ID NUMBER, USERNAME, PASSWORD, PROFILE IMAGE FILE,
# ELECTRICAL/WATER/GAS ACCOUNT PAYMENT, CAR FUEL TYPE, CAR TYPE, DISTANCE
$ws.Range("A1").Value = "FULL NAME"
$ws.Range("B1").Value = "ID NUMBER"
$ws.Range("C1").Value = "USERNAME"
$ws.Range("D1").Value = "PASSWORD"
$ws.Range("E1").Value = "PROFILE IMAGE FILE"
$ws.Range("I1").Value = "CAR FUEL TYPE"
$ws.Range("J1").Value = "CAR TYPE"
$ws.Range("K1").Value = "DISTANCE"
$ws.Range("H1").Value = "GAS ACCOUNT PAYMENT"
$ws.Range("G1").Value = "WATER ACCOUNT PAYMENT"
$ws.Range("F1").Value = "ELECTRICAL ACCOUNT PAYMENT"

# H1 picks up a tweaked alignment style in the new workbook
$ws.Range("H1").HorizontalAlignment = -4131

# Default font across the workbook switches from Arial to Calibri
$ws.Cells.Font.Name = "Calibri"
$ws.Range("D2").Font.Name = "Calibri"

# Selection moves to F4
$ws.Range("F4").Select()

Write-Output "done"
